$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "With the flight scanner API, an API that tracks flight data, the issue is two of the imports don’t work, Ive been searching everywhere on stack overflow, looking on the documentation, and reading, and still cant find a solution. One of the things I tried was deleting, and re-importing all of the dependencies of the project, However, when i ran the code it tells me other imports arent working"

# Add new row 48: date, hours, description
$ws.Cells.Item(48, 1).Value = 45239
$ws.Cells.Item(48, 1).NumberFormat = $ws.Cells.Item(47, 1).NumberFormat
$ws.Cells.Item(48, 2).Value = 3
$ws.Cells.Item(48, 3).Value = $newText

# Update the view: select C48 and scroll so A26 is the top-left visible cell
$ws.Range("C48").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 103

$wb.Save()
